$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the three shapes that make up the "Person State Machine" figure:
#   - the diagram picture
#   - the "Person State Machine" title textbox
#   - the "Residences / Entertainment" caption textbox
$pic = $null
$tbPerson = $null
$tbResidences = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Picture 17") {
        $pic = $shp
    } elseif ($shp.Name -eq "TextBox 21") {
        $tbPerson = $shp
    } elseif ($shp.Name -eq "TextBox 24") {
        $tbResidences = $shp
    }
}

# Enlarge the state-machine diagram picture, keeping its top-left corner fixed.
$pic.Width = 806.9595275590551
$pic.Height = 474.24236220472443

# Re-center the "Person State Machine" title over the now-wider picture.
$tbPerson.Left = 1951.4471653543308

# Push the caption textbox down so it clears the taller picture.
$tbResidences.Top = 1002.5508174015748

# Rename the shapes to reflect that they were effectively re-created/re-pasted
# (new shape identities), matching the authoring tool's behavior.
$pic.Name = "Picture 23"
$tbPerson.Name = "TextBox 26"
$tbResidences.Name = "TextBox 28"

# Move the picture and its two textboxes to the end of the z-order / shape
# tree, in the same relative order, so they render after every other shape
# on the slide (this is what happens when they are deleted and re-added).
$pic.ZOrder(0)
$tbPerson.ZOrder(0)
$tbResidences.ZOrder(0)
